$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (pushes all the existing weekly records
# down by one row, through row 29 -> row 30), carrying formatting down from
# the row above per Excel's default behavior.
$ws.Rows("12:12").Insert()

# Populate the newly-inserted row 12 with the new weekly record.
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").Value = 44914
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 100112003
$ws.Range("G12").Value = "Ajo"
$ws.Range("H12").Value = "Chino"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14500
$ws.Range("N12").Value = "$/caja 10 kilos"
$ws.Range("O12").Value = "China"
$ws.Range("P12").Value = 1450
$ws.Range("Q12").Value = 10
$ws.Range("R12").Value = "Hortaliza"
